$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: insert a brand new answer paragraph after the bullet question
# "Describir la diferencia entre Tarjeta de interfaz de red (NIC), puerto
# físico e interfaz de red." and before the next numbered question.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Describir la diferencia entre Tarjeta de interfaz de red (NIC), puerto físico e interfaz de red.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Anchor paragraph for change 1 not found" }

$para = $rng.Paragraphs(1)
$insertPoint = $para.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$newPara = $insertPoint.Paragraphs(1)

$newParaXml = '<w:p ' + $wns + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="NormalWeb"/>' +
        '<w:spacing w:before="91" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>' +
        '<w:ind w:left="1396" w:right="1170"/>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">La diferencia entre Tarjeta de Interfaz de Red, puerto físico e interfaz de red se da en que la primera es un adaptador LAN que proporciona conexión física con la red en la computadora u otro dispositivo host, mientras que el puerto físico es un conector en un dispositivo de red en el cual el medio se conecta con un host u otro dispositivo de red, finalmente la interfaz de red corresponde a puertos especializados de un dispositivo de </w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t>internetworking</w:t>' +
    '</w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve"> que se conecta con redes individuales. </w:t>' +
    '</w:r>' +
    '</w:p>'

$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# Change 2: split the "En cuanto a los modelos de referencia..." run into two
# runs, inserting a lastRenderedPageBreak before "abstracción mucho menor...".
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$searchText2 = "En cuanto a los modelos de referencia se tiene que proporcionar un nivel de abstracción mucho menor, generando así una referencia común para mantener consistencia en todos los tipos de protocolos y servicios de la red. Es decir que en estos modelos no se genera tanta minucia ni especificación, por lo tanto, no se puede definir una forma precisa de los servicios de la arquitectura de red. Un ejemplo de esto sería el modelo OSI."
$found2 = $rng2.Find.Execute($searchText2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Anchor paragraph for change 2 not found" }

$para2 = $rng2.Paragraphs(1)
$insertPoint2 = $para2.Range
$insertPoint2.Collapse(0)
$insertPoint2.InsertParagraphAfter()
$newPara2 = $insertPoint2.Paragraphs(1)

$part1 = "En cuanto a los modelos de referencia se tiene que proporcionar un nivel de "
$part2 = "abstracción mucho menor, generando así una referencia común para mantener consistencia en todos los tipos de protocolos y servicios de la red. Es decir que en estos modelos no se genera tanta minucia ni especificación, por lo tanto, no se puede definir una forma precisa de los servicios de la arquitectura de red. Un ejemplo de esto sería el modelo OSI."

$newPara2Xml = '<w:p ' + $wns + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:spacing w:before="91"/>' +
        '<w:ind w:left="1396" w:right="1170" w:firstLine="0"/>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">' + $part1 + '</w:t>' +
    '</w:r>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:lastRenderedPageBreak/>' +
        '<w:t xml:space="preserve">' + $part2 + '</w:t>' +
    '</w:r>' +
    '</w:p>'

$newPara2.Range.InsertXML($newPara2Xml)
$para2.Range.Delete()

# ---------------------------------------------------------------------------
# Change 3: remove the lastRenderedPageBreak from the bullet paragraph
# "Describir la función de cada capa en los dos modelos de red: TCP/IP y OSI."
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$searchText3 = "Describir la función de cada capa en los dos modelos de red: TCP/IP y OSI."
$found3 = $rng3.Find.Execute($searchText3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Anchor paragraph for change 3 not found" }

$para3 = $rng3.Paragraphs(1)
$insertPoint3 = $para3.Range
$insertPoint3.Collapse(0)
$insertPoint3.InsertParagraphAfter()
$newPara3 = $insertPoint3.Paragraphs(1)

$newPara3Xml = '<w:p ' + $wns + '>' +
    '<w:pPr>' +
        '<w:numPr>' +
            '<w:ilvl w:val="0"/>' +
            '<w:numId w:val="3"/>' +
        '</w:numPr>' +
        '<w:spacing w:before="91"/>' +
        '<w:ind w:right="1170"/>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t>' + $searchText3 + '</w:t>' +
    '</w:r>' +
    '</w:p>'

$newPara3.Range.InsertXML($newPara3Xml)
$para3.Range.Delete()

# ---------------------------------------------------------------------------
# Change 4: add a lastRenderedPageBreak before "1. Física:" (keeping the
# second run of that paragraph untouched).
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$searchText4 = "1. Física: en esta capa se describen los medios físicos para activar, mantener y desactivar conexiones por medio de las cuales se va a realizar la transmisión de bits hacia y desde un dispositivo de red."
$found4 = $rng4.Find.Execute($searchText4, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "Anchor paragraph for change 4 not found" }

$para4 = $rng4.Paragraphs(1)
$insertPoint4 = $para4.Range
$insertPoint4.Collapse(0)
$insertPoint4.InsertParagraphAfter()
$newPara4 = $insertPoint4.Paragraphs(1)

$newPara4Xml = '<w:p ' + $wns + '>' +
    '<w:pPr>' +
        '<w:spacing w:before="91"/>' +
        '<w:ind w:left="1440" w:right="1170"/>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:b/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:lastRenderedPageBreak/>' +
        '<w:t>1. Física:</w:t>' +
    '</w:r>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve"> en esta capa se describen los medios físicos para activar, mantener y desactivar conexiones por medio de las cuales se va a realizar la transmisión de bits hacia y desde un dispositivo de red.</w:t>' +
    '</w:r>' +
    '</w:p>'

$newPara4.Range.InsertXML($newPara4Xml)
$para4.Range.Delete()

Write-Host "All changes applied"
